$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "CghVf808"
$ws.Range("B2").Value = 231103130
$ws.Range("C2").Value = "sgbfgjw24"
$ws.Range("D2").Value = "u!W%5P4v"
$ws.Range("F2").Value = "jBgVYOnE"
$ws.Range("G2").Value = "YMtq"

# Row 3 updates
$ws.Range("A3").Value = "zwnQT970"
$ws.Range("B3").Value = 231103129
$ws.Range("C3").Value = "ytetdcn27"
$ws.Range("D3").Value = "X#ue67E!"
$ws.Range("F3").Value = "ioqSJOfK"
$ws.Range("G3").Value = "rVoh"

# Row 4 updates
$ws.Range("A4").Value = "dbunr222"
$ws.Range("B4").Value = 231103128
$ws.Range("C4").Value = "rqzogkd85"
$ws.Range("D4").Value = "f!%9Dq4A"
$ws.Range("F4").Value = "dmnFKVTh"
$ws.Range("G4").Value = "jdNl"
